$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (header "Förändrad") holds a date serial that was bumped by one day
# (46081 -> 46082) for every data row (rows 2 through 393).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 393
}

$ws.Range("C2:C$lastRow").Value = 46082
